$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded. Insert a row above the
# current row 14 (shifting the existing rows 14-36 down to 15-37) and
# fill it in with the new reading, matching the columns used by every
# other row in this block.
$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14, 1).Value = 11
$ws.Cells.Item(14, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(14, 3).Value = "Bíobío"
$ws.Cells.Item(14, 4).Value = 44775
$ws.Cells.Item(14, 5).Value = 8
$ws.Cells.Item(14, 6).Value = 100114007
$ws.Cells.Item(14, 7).Value = "Jengibre"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 20
$ws.Cells.Item(14, 11).Value = 12000
$ws.Cells.Item(14, 12).Value = 13000
$ws.Cells.Item(14, 13).Value = 12500
$ws.Cells.Item(14, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(14, 15).Value = "Perú"
$ws.Cells.Item(14, 16).Value = 962
$ws.Cells.Item(14, 17).Value = 13
$ws.Cells.Item(14, 18).Value = "Hortaliza"
